$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new columns: apoio_std/min/max after L (apoio_medio), and contribuicoes_std/min/max after contribuicoes_med
$ws.Columns("M:O").Insert()
$ws.Columns("R:T").Insert()

# --- Header renames / new headers ---
$ws.Range("H1").Value = "arrecadado_avg"
$ws.Range("I1").Value = "arrecadado_std"
$ws.Range("J1").Value = "arrecadado_min"
$ws.Range("K1").Value = "arrecadado_max"
$ws.Range("M1").Value = "apoio_std"
$ws.Range("N1").Value = "apoio_min"
$ws.Range("O1").Value = "apoio_max"
$ws.Range("Q1").Value = "contribuicoes_med"
$ws.Range("R1").Value = "contribuicoes_std"
$ws.Range("S1").Value = "contribuicoes_min"
$ws.Range("T1").Value = "contribuicoes_max"

# --- Data rows: apoio_medio changes, apoio_std/min/max + contribuicoes_std/min/max are new ---
$r = 2
$ws.Cells.Item($r, 12).Value = 49.34563627103121
$ws.Cells.Item($r, 13).Value = 3.157738703897078
$ws.Cells.Item($r, 14).Value = 45.33997822063326
$ws.Cells.Item($r, 15).Value = 52.78244028225378
$ws.Cells.Item($r, 18).Value = 61.60925255186918
$ws.Cells.Item($r, 19).Value = 1
$ws.Cells.Item($r, 20).Value = 171

$r = 3
$ws.Cells.Item($r, 12).Value = 75.92822821742841
$ws.Cells.Item($r, 13).Value = 33.04151967611716
$ws.Cells.Item($r, 14).Value = 44.41698989306864
$ws.Cells.Item($r, 15).Value = 110.3124119232409
$ws.Cells.Item($r, 18).Value = 19.65536398374076
$ws.Cells.Item($r, 19).Value = 14
$ws.Cells.Item($r, 20).Value = 49

$r = 4
$ws.Cells.Item($r, 12).Value = 67.78664768952572
$ws.Cells.Item($r, 13).Value = 31.70775623287225
$ws.Cells.Item($r, 14).Value = 23.84123780968962
$ws.Cells.Item($r, 15).Value = 139.5126743984584
$ws.Cells.Item($r, 18).Value = 83.77267994257976
$ws.Cells.Item($r, 19).Value = 1
$ws.Cells.Item($r, 20).Value = 366

$r = 5
$ws.Cells.Item($r, 12).Value = 69.40537068990857
$ws.Cells.Item($r, 13).Value = 29.60589045154097
$ws.Cells.Item($r, 14).Value = 27.45405411991957
$ws.Cells.Item($r, 15).Value = 157.4361495951104
$ws.Cells.Item($r, 18).Value = 134.3866279276542
$ws.Cells.Item($r, 19).Value = 1
$ws.Cells.Item($r, 20).Value = 575

$r = 6
$ws.Cells.Item($r, 12).Value = 70.98726579114015
$ws.Cells.Item($r, 13).Value = 26.1115491044645
$ws.Cells.Item($r, 14).Value = 11.93343625774652
$ws.Cells.Item($r, 15).Value = 119.480863051166
$ws.Cells.Item($r, 18).Value = 87.9202644809713
$ws.Cells.Item($r, 19).Value = 1
$ws.Cells.Item($r, 20).Value = 303

$r = 7
$ws.Cells.Item($r, 12).Value = 57.05344062495174
$ws.Cells.Item($r, 13).Value = 9.218686058143824
$ws.Cells.Item($r, 14).Value = 48.67908690791425
$ws.Cells.Item($r, 15).Value = 66.93155962911479
$ws.Cells.Item($r, 18).Value = 48.80915214724933
$ws.Cells.Item($r, 19).Value = 58
$ws.Cells.Item($r, 20).Value = 155

$r = 8
$ws.Cells.Item($r, 12).Value = 84.3699140329713
$ws.Cells.Item($r, 13).Value = 38.84510413047353
$ws.Cells.Item($r, 14).Value = 29.26216513679551
$ws.Cells.Item($r, 15).Value = 121.4283973998111
$ws.Cells.Item($r, 18).Value = 32.65424934062947
$ws.Cells.Item($r, 19).Value = 16
$ws.Cells.Item($r, 20).Value = 96

$r = 9
$ws.Cells.Item($r, 12).Value = 48.84608159861521
$ws.Cells.Item($r, 13).Value = 6.201395329160901
$ws.Cells.Item($r, 14).Value = 43.30757970997428
$ws.Cells.Item($r, 15).Value = 56.35300160617668
$ws.Cells.Item($r, 18).Value = 19.27001124372618
$ws.Cells.Item($r, 19).Value = 32
$ws.Cells.Item($r, 20).Value = 73

$r = 10
$ws.Cells.Item($r, 12).Value = 69.83870009492911
$ws.Cells.Item($r, 13).Value = 35.4080170462333
$ws.Cells.Item($r, 14).Value = 12.19662302883409
$ws.Cells.Item($r, 15).Value = 196.4212117364618
$ws.Cells.Item($r, 18).Value = 108.5202951295182
$ws.Cells.Item($r, 19).Value = 1
$ws.Cells.Item($r, 20).Value = 571

$r = 11
$ws.Cells.Item($r, 12).Value = 65.64033709136109
$ws.Cells.Item($r, 13).Value = 14.83568373171277
$ws.Cells.Item($r, 14).Value = 52.25756349732896
$ws.Cells.Item($r, 15).Value = 91.02088659474175
$ws.Cells.Item($r, 18).Value = 47.15612367444974
$ws.Cells.Item($r, 19).Value = 10
$ws.Cells.Item($r, 20).Value = 141

$r = 12
$ws.Cells.Item($r, 12).Value = 77.84263902306331
$ws.Cells.Item($r, 13).Value = 5.602932423344289
$ws.Cells.Item($r, 14).Value = 73.88076751198659
$ws.Cells.Item($r, 15).Value = 81.80451053414004
$ws.Cells.Item($r, 18).Value = 139.3000358937499
$ws.Cells.Item($r, 19).Value = 30
$ws.Cells.Item($r, 20).Value = 227

$r = 13
$ws.Cells.Item($r, 12).Value = 45.35433730016975
$ws.Cells.Item($r, 13).Value = 25.46833055396829
$ws.Cells.Item($r, 14).Value = 20.15182714817413
$ws.Cells.Item($r, 15).Value = 83.50597686157313
$ws.Cells.Item($r, 18).Value = 102.1846368100411
$ws.Cells.Item($r, 19).Value = 5
$ws.Cells.Item($r, 20).Value = 236

$r = 14
$ws.Cells.Item($r, 12).Value = 51.46088164215199
$ws.Cells.Item($r, 13).Value = 23.24567672487991
$ws.Cells.Item($r, 14).Value = 13.05503087794559
$ws.Cells.Item($r, 15).Value = 111.3662854612667
$ws.Cells.Item($r, 18).Value = 117.5574611271725
$ws.Cells.Item($r, 19).Value = 2
$ws.Cells.Item($r, 20).Value = 539

$r = 15
$ws.Cells.Item($r, 12).Value = 61.5562500780191
$ws.Cells.Item($r, 13).Value = 23.34464831810631
$ws.Cells.Item($r, 14).Value = 16.18065842403185
$ws.Cells.Item($r, 15).Value = 138.0184648379384
$ws.Cells.Item($r, 18).Value = 68.49560506255611
$ws.Cells.Item($r, 19).Value = 1
$ws.Cells.Item($r, 20).Value = 328

$r = 16
$ws.Cells.Item($r, 12).Value = 59.77015408752801
$ws.Cells.Item($r, 13).Value = 14.08491272336895
$ws.Cells.Item($r, 14).Value = 45.90673655161051
$ws.Cells.Item($r, 15).Value = 82.15413766761272
$ws.Cells.Item($r, 18).Value = 100.0674772341144
$ws.Cells.Item($r, 19).Value = 10
$ws.Cells.Item($r, 20).Value = 284

$r = 17
$ws.Cells.Item($r, 12).Value = 82.66311410722277
$ws.Cells.Item($r, 13).Value = 41.59443116618839
$ws.Cells.Item($r, 14).Value = 20.33774597757668
$ws.Cells.Item($r, 15).Value = 201.2220224911509
$ws.Cells.Item($r, 18).Value = 223.4736944471927
$ws.Cells.Item($r, 19).Value = 1
$ws.Cells.Item($r, 20).Value = 1318

$r = 18
$ws.Cells.Item($r, 12).Value = 78.50979961663667
$ws.Cells.Item($r, 13).Value = 35.09459453111103
$ws.Cells.Item($r, 14).Value = 10.77163914429046
$ws.Cells.Item($r, 15).Value = 233.3973531230909
$ws.Cells.Item($r, 18).Value = 239.8765062590198
$ws.Cells.Item($r, 19).Value = 1
$ws.Cells.Item($r, 20).Value = 2120

$r = 19
$ws.Cells.Item($r, 12).Value = 82.28079522919644
$ws.Cells.Item($r, 13).Value = 10.15766237346972
$ws.Cells.Item($r, 14).Value = 74.11921242291902
$ws.Cells.Item($r, 15).Value = 93.65675828662884
$ws.Cells.Item($r, 18).Value = 118.1199954848176
$ws.Cells.Item($r, 19).Value = 2
$ws.Cells.Item($r, 20).Value = 226

$r = 20
$ws.Cells.Item($r, 12).Value = 41.02094512827082
$ws.Cells.Item($r, 13).Value = 18.86799504524634
$ws.Cells.Item($r, 14).Value = 21.99487001791516
$ws.Cells.Item($r, 15).Value = 65.85076384224313
$ws.Cells.Item($r, 18).Value = 29.91515780781821
$ws.Cells.Item($r, 19).Value = 2
$ws.Cells.Item($r, 20).Value = 75

$r = 21
$ws.Cells.Item($r, 12).Value = 81.7894238310976
$ws.Cells.Item($r, 13).Value = 45.73065301198442
$ws.Cells.Item($r, 14).Value = 18.89233795141325
$ws.Cells.Item($r, 15).Value = 306.2848483936102
$ws.Cells.Item($r, 18).Value = 153.2010945781367
$ws.Cells.Item($r, 19).Value = 1
$ws.Cells.Item($r, 20).Value = 1004

$r = 22
$ws.Cells.Item($r, 12).Value = 85.51665070147128
$ws.Cells.Item($r, 13).Value = 52.63790605978169
$ws.Cells.Item($r, 14).Value = 21.00493274015408
$ws.Cells.Item($r, 15).Value = 254.2443749773306
$ws.Cells.Item($r, 18).Value = 98.6216776074999
$ws.Cells.Item($r, 19).Value = 1
$ws.Cells.Item($r, 20).Value = 344

$r = 23
$ws.Cells.Item($r, 12).Value = 36.90833212357519
$ws.Cells.Item($r, 13).Value = 0
$ws.Cells.Item($r, 14).Value = 36.90833212357519
$ws.Cells.Item($r, 15).Value = 36.90833212357519
$ws.Cells.Item($r, 18).Value = 0
$ws.Cells.Item($r, 19).Value = 55
$ws.Cells.Item($r, 20).Value = 55

$r = 24
$ws.Cells.Item($r, 12).Value = 80.55678320085916
$ws.Cells.Item($r, 13).Value = 41.26322575364973
$ws.Cells.Item($r, 14).Value = 11.52676430516467
$ws.Cells.Item($r, 15).Value = 461.5197709071476
$ws.Cells.Item($r, 18).Value = 419.7037665583493
$ws.Cells.Item($r, 19).Value = 1
$ws.Cells.Item($r, 20).Value = 7954
